# This script reproduces the data refresh captured by the commit message
# "Updated cryptos list ... with GitHub Actions". It rewrites the Price
# (column D) and Volume(1h) (column E) text for the rows whose figures
# changed, and, for the two row pairs whose coin/link/price/volume data
# moved to the neighboring row (TRON <-> WrappedEther at rows 12-13, and
# TheSandbox <-> VeChain at rows 39-40), it also rewrites Coin (B) and
# Link (C).
#
# Column D/E values are stored as plain text in the workbook (e.g.
# "26.977.89" or "  -0.26%  "), not real numbers. A leading apostrophe is
# included in every assignment below so Excel always stores the value as
# literal text instead of silently parsing number-looking strings into a
# Double (which would also corrupt values like "20.79" via floating point
# rounding). The apostrophe itself is never part of the stored text. After
# the value is set, the cell style is reset to "Normal" so the forced-text
# formatting Excel applies doesn't leave a stray style on the cell (cells
# in the source file carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.977.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.26%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.874.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.55%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.19%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'305.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.32%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.09%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5091"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3666"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.38%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07216"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8959"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.84%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.58%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'1.886.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.26%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'TRON"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.07524"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.97%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'95.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.36%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.10%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.20%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008540"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.96%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'14.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.11%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'27.009.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.26%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.028"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.17%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.130.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.67%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.27%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.405"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.76%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.91%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.792"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.96%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.47%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.097"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.65%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.61%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.730"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.57%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.743"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.80%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.09173"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.66%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.48%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7505"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.36%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.971"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.20%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.161"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.54%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D38").Value = "'2.541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.83%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'VeChain"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.02004"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.76%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.5620"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.35%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.079"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.31%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'6.648"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.44%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'115.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.92%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.586"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.84%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.1481"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.60%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4769"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.97%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.0000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.15%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'10.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.01%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.572"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'36.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.13%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'63.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.00%  "
$ws.Range("E51").Style = "Normal"
